{"js": "// Replace each old equation string with its updated result, matching the\n// diff exactly (one <w:t> run per multiplication problem in the table).\nconst replacements = [\n  [\"797\u00d77=5579\", \"302\u00d72=604\"],\n  [\"714\u00d73=2142\", \"599\u00d79=5391\"],\n  [\"513\u00d73=1539\", \"323\u00d75=1615\"],\n  [\"196\u00d75=980\", \"851\u00d77=5957\"],\n  [\"900\u00d76=5400\", \"190\u00d77=1330\"],\n  [\"899\u00d77=6293\", \"239\u00d72=478\"],\n  [\"200\u00d75=1000\", \"917\u00d78=7336\"],\n  [\"712\u00d77=4984\", \"678\u00d74=2712\"],\n  [\"500\u00d72=1000\", \"499\u00d76=2994\"],\n  [\"951\u00d78=7608\", \"249\u00d75=1245\"],\n  [\"367\u00d78=2936\", \"784\u00d76=4704\"],\n  [\"844\u00d77=5908\", \"224\u00d72=448\"],\n  [\"529\u00d75=2645\", \"903\u00d73=2709\"],\n  [\"148\u00d78=1184\", \"269\u00d77=1883\"],\n  [\"303\u00d79=2727\", \"974\u00d76=5844\"],\n  [\"275\u00d77=1925\", \"189\u00d75=945\"],\n  [\"840\u00d75=4200\", \"470\u00d79=4230\"],\n  [\"315\u00d75=1575\", \"194\u00d74=776\"],\n  [\"673\u00d73=2019\", \"970\u00d79=8730\"],\n  [\"287\u00d72=574\", \"186\u00d78=1488\"],\n  [\"185\u00d75=925\", \"581\u00d73=1743\"],\n  [\"636\u00d79=5724\", \"666\u00d76=3996\"],\n  [\"200\u00d78=1600\", \"391\u00d75=1955\"],\n  [\"325\u00d76=1950\", \"716\u00d73=2148\"],\n  [\"429\u00d72=858\", \"717\u00d76=4302\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each three-digit-by-one-digit multiplication problem in the\n# answer table to the newly generated values (1:1 with the diff, matched\n# on the previous full \"A\u00d7B=C\" string so each cell is targeted uniquely).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"797\u00d77=5579\", \"302\u00d72=604\"),\n    @(\"714\u00d73=2142\", \"599\u00d79=5391\"),\n    @(\"513\u00d73=1539\", \"323\u00d75=1615\"),\n    @(\"196\u00d75=980\", \"851\u00d77=5957\"),\n    @(\"900\u00d76=5400\", \"190\u00d77=1330\"),\n    @(\"899\u00d77=6293\", \"239\u00d72=478\"),\n    @(\"200\u00d75=1000\", \"917\u00d78=7336\"),\n    @(\"712\u00d77=4984\", \"678\u00d74=2712\"),\n    @(\"500\u00d72=1000\", \"499\u00d76=2994\"),\n    @(\"951\u00d78=7608\", \"249\u00d75=1245\"),\n    @(\"367\u00d78=2936\", \"784\u00d76=4704\"),\n    @(\"844\u00d77=5908\", \"224\u00d72=448\"),\n    @(\"529\u00d75=2645\", \"903\u00d73=2709\"),\n    @(\"148\u00d78=1184\", \"269\u00d77=1883\"),\n    @(\"303\u00d79=2727\", \"974\u00d76=5844\"),\n    @(\"275\u00d77=1925\", \"189\u00d75=945\"),\n    @(\"840\u00d75=4200\", \"470\u00d79=4230\"),\n    @(\"315\u00d75=1575\", \"194\u00d74=776\"),\n    @(\"673\u00d73=2019\", \"970\u00d79=8730\"),\n    @(\"287\u00d72=574\", \"186\u00d78=1488\"),\n    @(\"185\u00d75=925\", \"581\u00d73=1743\"),\n    @(\"636\u00d79=5724\", \"666\u00d76=3996\"),\n    @(\"200\u00d78=1600\", \"391\u00d75=1955\"),\n    @(\"325\u00d76=1950\", \"716\u00d73=2148\"),\n    @(\"429\u00d72=858\", \"717\u00d76=4302\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"No match found for `\"$oldText`\"\"\n    }\n}\n"}
